# Commit: "Add new column 'Servised by' to Card24 by admin"
#
# On the Card24 sheet the data rows (2-12) had the "Correction" values
# (column N) sitting empty while the "Servised by" values (column O) held
# the placeholder "nan" text. The edit moves the placeholder text from the
# "Servised by" column into the "Correction" column and leaves the new
# "Servised by" column blank for each data row - the header row (row 1) is
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 14).Value = "nan"   # column N - "Correction"
    $ws.Cells.Item($row, 15).ClearContents() # column O - "Servised by"
}
